# Apply cryptocurrency price/volume updates to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.212.40'
$ws.Range("E2").Value = '  -0.44%  '

$ws.Range("D3").Value = '2.388.31'
$ws.Range("E3").Value = '  -3.67%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.25'
$ws.Range("E5").Value = '  -0.69%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.60'
$ws.Range("E6").Value = '  -3.39%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D9").Value = '2.386.74'
$ws.Range("E9").Value = '  -3.69%  '

$ws.Range("E10").Value = '  -2.34%  '

$ws.Range("E11").Value = '  +0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.29'
$ws.Range("E12").Value = '  -3.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.346'
$ws.Range("E13").Value = '  -3.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.48'
$ws.Range("E14").Value = '  -2.76%  '

$ws.Range("D15").Value = '2.818.57'
$ws.Range("E15").Value = '  -3.68%  '

$ws.Range("E16").Value = '  -1.48%  '

$ws.Range("D17").Value = '61.167.37'
$ws.Range("E17").Value = '  -0.35%  '

$ws.Range("D18").Value = '2.387.02'
$ws.Range("E18").Value = '  -3.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.77'
$ws.Range("E19").Value = '  -3.67%  '

$ws.Range("E20").Value = '  -1.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '318.74'
$ws.Range("E21").Value = '  -1.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.71'
$ws.Range("E22").Value = '  -5.03%  '

$ws.Range("E24").Value = '  +2.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.53'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.23'
$ws.Range("E26").Value = '  +5.45%  '

$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("D28").Value = '2.505.58'
$ws.Range("E28").Value = '  -3.91%  '

$ws.Range("D29").Value = '0.0₃0929'
$ws.Range("E29").Value = '  -6.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '525.10'
$ws.Range("E30").Value = '  -3.59%  '

$ws.Range("E31").Value = '  -5.49%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.09'
$ws.Range("E32").Value = '  -3.08%  '

$ws.Range("E33").Value = '  -4.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.84'
$ws.Range("E34").Value = '  -3.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").Value = '  -1.11%  '

$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.52'
$ws.Range("E37").Value = '  -6.37%  '

$ws.Range("E38").Value = '  -4.11%  '

$ws.Range("E39").Value = '  -1.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.84'
$ws.Range("E40").Value = '  +6.00%  '

$ws.Range("E41").Value = '  -2.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '139.14'
$ws.Range("E42").Value = '  -5.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.27'
$ws.Range("E44").Value = '  -0.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.16'
$ws.Range("E45").Value = '  -9.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '140.80'
$ws.Range("E46").Value = '  -4.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.62'
$ws.Range("E47").Value = '  -0.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.16'
$ws.Range("E48").Value = '  -5.28%  '

$ws.Range("E49").Value = '  -3.92%  '

$ws.Range("E50").Value = '  -3.74%  '

$ws.Range("E51").Value = '  -1.04%  '

# Restore default (unstyled) appearance for cells where we had to
# temporarily force a text number format to stop Excel from
# reinterpreting numeric-looking strings (e.g. "141.60") as numbers.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
